$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as literal text in the source data
# (e.g. thousands-separated "43.724.63", or trailing-zero "8.00"). Assigning
# such strings straight to Range.Value lets Excel auto-coerce anything that
# parses as a number, silently dropping trailing zeros / changing content.
# Force the column to Text ("@") before writing, then clear the temporary
# formatting back off (restores the original unstyled cells) once done.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.724.63"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.293.38"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "103.56"
$ws.Range("E5").Value = "  +5.89%  "
$ws.Range("D6").Value = "270.77"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "45.68"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "8.00"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "15.70"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "0.859"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "2.290.00"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "43.773.45"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "72.31"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "2.49"
$ws.Range("E21").Value = "  +10.24%  "
$ws.Range("D22").Value = "233.61"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "2.90"
$ws.Range("E23").Value = "  +14.69%  "
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "11.22"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "39.90"
$ws.Range("E28").Value = "  +5.02%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "177.10"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Value = "21.83"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "0.0901"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "5.47"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  +13.23%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "3.56"
$ws.Range("E38").Value = "  +5.93%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.237"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "12.27"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "65.60"
$ws.Range("E43").Value = "  +5.71%  "
$ws.Range("D44").Value = "8.82"
$ws.Range("D45").Value = "5.25"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "1.23"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").Value = "98.70"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  +8.73%  "
$ws.Range("D50").Value = "1.53"
$ws.Range("E50").Value = "  +11.17%  "
$ws.Range("D51").Value = "2.516.98"
$ws.Range("E51").Value = "  -0.76%  "

$dRange.ClearFormats()
